$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 618.8099999999999
$ws.Range("C3").Value = 643.89
$ws.Range("C4").Value = 612.58
$ws.Range("C5").Value = 628.23
$ws.Range("C6").Value = 628.23
